# Estadisticos Segundo Parcial 23 Mayo
# Updates:
#  - "Estadisticos 2P": fill in the second-partial statistics (previously all zero/default)
#  - "Estadisticos Final": refresh combined stats for most groups
#  - "Rescatables": replace the make-up-exam candidate list with the new roster

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P" — second partial statistics now available
# ---------------------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

# row 2 - Ingles II / 2ARHM
$ws2P.Cells.Item(2,4).Value = 0
$ws2P.Cells.Item(2,5).Value = 9
$ws2P.Cells.Item(2,6).Value = 40
$ws2P.Cells.Item(2,7).Value = 81.63
$ws2P.Cells.Item(2,8).Value = 7

# row 3 - Ingles II / 2BLCM
$ws2P.Cells.Item(3,4).Value = 0
$ws2P.Cells.Item(3,5).Value = 4
$ws2P.Cells.Item(3,6).Value = 34
$ws2P.Cells.Item(3,7).Value = 89.47
$ws2P.Cells.Item(3,8).Value = 7.8

# row 4 - Ingles IV / 4ARHM
$ws2P.Cells.Item(4,4).Value = 0
$ws2P.Cells.Item(4,5).Value = 1
$ws2P.Cells.Item(4,6).Value = 39
$ws2P.Cells.Item(4,7).Value = 97.5
$ws2P.Cells.Item(4,8).Value = 9.6

# row 5 - Ingles IV / 4BLCM
$ws2P.Cells.Item(5,4).Value = 0
$ws2P.Cells.Item(5,5).Value = 1
$ws2P.Cells.Item(5,6).Value = 37
$ws2P.Cells.Item(5,7).Value = 97.37
$ws2P.Cells.Item(5,8).Value = 8.8

# ---------------------------------------------------------------------------
# Sheet "Estadisticos Final" — recompute with the new second-partial figures
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

# row 2 - Ingles II / 2ARHM
$wsFinal.Cells.Item(2,5).Value = 9
$wsFinal.Cells.Item(2,6).Value = 40
$wsFinal.Cells.Item(2,7).Value = 81.63
$wsFinal.Cells.Item(2,8).Value = 7.7

# row 3 - Ingles II / 2BLCM (only the Promedio changes)
$wsFinal.Cells.Item(3,8).Value = 8.2

# row 4 - Ingles IV / 4ARHM -> unchanged

# row 5 - Ingles IV / 4BLCM
$wsFinal.Cells.Item(5,5).Value = 1
$wsFinal.Cells.Item(5,6).Value = 37
$wsFinal.Cells.Item(5,7).Value = 97.37
$wsFinal.Cells.Item(5,8).Value = 9.2

# ---------------------------------------------------------------------------
# Sheet "Rescatables" — updated list of students eligible for make-up exams
# ---------------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

$rescatables = @(
    @(24330051920117, "PLIEGO",   "TORRES",  "MYA YAMILET",    "Ingles II", "2ARHM", 4),
    @(24330051920348, "TRUJILLO", "FLORES",  "VALERIA",        "Ingles II", "2ARHM", 4),
    @(24330051920341, "VALDIVIA", "TENORIO", "IKER YHAIR",     "Ingles II", "2BLCM", 4),
    @(24330051920345, "GARIBAY",  "GOMEZ",   "LIZBETH MARIAM", "Ingles II", "2ARHM", 3),
    @(23330051920263, "OLMOS",    "ORTEGA",  "ANGEL GABRIEL",  "Ingles IV", "4ARHM", 3),
    @(24330051920169, "BAEZ",     "LOPEZ",   "ULISES EZEQUIEL","Ingles II", "2ARHM", 2),
    @(24330051920350, "FLORES",   "LOBATO",  "MARIANA",        "Ingles II", "2ARHM", 2),
    @(24330051920132, "GONZALEZ", "CRUZ",    "JESUS",          "Ingles II", "2ARHM", 2),
    @(24330051920351, "PLIEGO",   "LORENZO", "CALEB SANTIAGO", "Ingles II", "2ARHM", 2),
    @(24330051920120, "ROSAS",    "GARRIDO", "YOSELYN",        "Ingles II", "2ARHM", 1)
)

$r = 2
foreach ($row in $rescatables) {
    $wsResc.Cells.Item($r,1).Value = $row[0]
    $wsResc.Cells.Item($r,2).Value = $row[1]
    $wsResc.Cells.Item($r,3).Value = $row[2]
    $wsResc.Cells.Item($r,4).Value = $row[3]
    $wsResc.Cells.Item($r,5).Value = $row[4]
    $wsResc.Cells.Item($r,6).Value = $row[5]
    $wsResc.Cells.Item($r,7).Value = $row[6]
    $r = $r + 1
}
